# Readme/controller_ubuntu.xlsx
# "Add new type of build and remove Minnow3."
#
# The sheet currently lists 4 board rows:
#   2: Minnow3 Module / FAB A
#   3: Minnow3 Module / FAB C
#   4: Leaf Hill      / FAB D   (Debug column = "N/A")
#   5: UP2            / FAB A
#
# Target: remove both "Minnow3 Module" rows entirely, keep "Leaf Hill" and
# "UP2" (which shift up into rows 2-3), and change the now-promoted "Leaf
# Hill" row's Debug column from "N/A" to "Y" (matching the rest of the
# table). The trailing blank spacer row and the legend row shift up with
# everything else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two "Minnow3 Module" rows (rows 2 and 3). Excel shifts rows
# 4 (Leaf Hill/FAB D) and 5 (UP2/FAB A) up into their place, along with
# the blank spacer row and the legend row below them.
$ws.Rows("2:3").Delete() | Out-Null

# The promoted "Leaf Hill" row (now row 2) had "N/A" in the Debug column
# with its own distinct shading. Bring it in line with every other row's
# "Y" / shared formatting: copy the format from the adjacent "Release"
# cell (already styled "Y") and then set the value.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = "Y"

# Restore the selection to the top data row (matches the saved view).
$ws.Range("D2").Select() | Out-Null
